# dlgCamperInsurance_pagProductData.xlsx - part 1 of AppiumLibrary -> python-appium-client swap.
# The only functional content change is the start-date seed value in B5, which moves from a
# hard-coded date to a templated "+32 days from today" placeholder used by the new test runner.
# Everything else here (column widths, picture sizing, selection) are the natural knock-on
# effects of that text getting wider, reproduced through the same user actions (resize the
# columns to fit, which nudges the two anchored screenshots, then the cursor lands on C7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Replace the static start date with the dynamic "+32 days" token used by the new framework.
$ws.Range("B5").Value = "<TODAY +32,+0,+0,'%m/%d/%Y'>"

# 2) Widen columns B and E so the new (longer) values keep fitting, matching the authored widths.
$ws.Columns.Item(2).ColumnWidth = 27.333333333333336
$ws.Columns.Item(5).ColumnWidth = 30.0

# 3) The two anchored screenshots are pinned "move but don't size with cells" (editAs=oneCell),
#    so widening columns B/E (which fall inside their cell spans) shifts their right/bottom edge.
#    Re-assert the shapes' rendered width so the stored anchor reflects the new column metrics.
$shp1 = $ws.Shapes.Item(1)
$shp1.Width = 814.7867187500001
$shp2 = $ws.Shapes.Item(2)
$shp2.Width = 555.7793762303149

# 4) Leave the selection where the author's cursor ended up after the edit.
$ws.Range("C7").Select() | Out-Null
